$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 / column E: the date value (05/14/2023, stored as serial 45255) is
# replaced by the literal text "05/14/2021" - simulating a faulty/raw value
# that should be surfaced in the upload error dialog instead of a parsed date.
# Force a text number format first so Excel keeps it as a string instead of
# re-parsing "05/14/2021" back into a date serial.
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "05/14/2021"

# Row 19 / column H: the time value (2.041666666666667) is replaced by the
# literal text "TIME" - another faulty value shown verbatim in the error
# dialog. Clear the existing time number format/style so the cell goes back
# to the default (unstyled) text cell, matching a freshly entered string.
$ws.Range("H19").ClearFormats()
$ws.Range("H19").Value = "TIME"

# Reflect the current selection/active cell as left after the edit.
$ws.Range("H20").Select()
